$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.218.04'
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').Value = '1.860.50'
$ws.Range('E3').Value = '  +0.67%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '0.7115'
$ws.Range('E5').Value = '  +1.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range('D6').Value = '237.86'
$ws.Range('E6').Value = '  -0.30%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range('D7').Value = '0.9999'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range('D8').Value = '0.08161'
$ws.Range('E8').Value = '  +9.65%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range('D9').Value = '0.3046'
$ws.Range('E9').Value = '  -0.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range('D10').Value = '23.20'
$ws.Range('E10').Value = '  -0.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range('D11').Value = '0.08170'
$ws.Range('E11').Value = '  +0.38%  '
$ws.Range('D12').Value = '1.902.43'
$ws.Range('E12').Value = '  +3.30%  '
$ws.Range('E13').Value = '  -0.82%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range('D14').Value = '0.7077'
$ws.Range('E14').Value = '  -2.74%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range('D15').Value = '89.63'
$ws.Range('E15').Value = '  +0.76%  '
$ws.Range('D16').Value = '29.201.24'
$ws.Range('E16').Value = '  +0.66%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range('D17').Value = '0.000007914'
$ws.Range('E17').Value = '  +3.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range('D18').Value = '5.793'
$ws.Range('E18').Value = '  +0.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range('D19').Value = '13.36'
$ws.Range('E19').Value = '  +2.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range('D20').Value = '237.31'
$ws.Range('E20').Value = '  -0.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range('D21').Value = '1.0000'
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').Value = '2.111.33'
$ws.Range('E22').Value = '  +2.02%  '
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range('D24').Value = '7.427'
$ws.Range('E24').Value = '  -2.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range('D25').Value = '162.43'
$ws.Range('E25').Value = '  +1.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range('D26').Value = '0.1465'
$ws.Range('E26').Value = '  +0.85%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range('D27').Value = '8.960'
$ws.Range('E27').Value = '  -0.59%  '
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range('D29').Value = '1.958'
$ws.Range('E29').Value = '  -0.70%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range('D30').Value = '1.426'
$ws.Range('E30').Value = '  +1.59%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range('D31').Value = '1.482'
$ws.Range('E31').Value = '  -0.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range('D32').Value = '4.399'
$ws.Range('E32').Value = '  -2.77%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range('D33').Value = '4.023'
$ws.Range('E33').Value = '  +0.44%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range('D34').Value = '0.05228'
$ws.Range('E34').Value = '  +0.64%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range('D35').Value = '1.169'
$ws.Range('E35').Value = '  -1.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range('D36').Value = '0.7085'
$ws.Range('E36').Value = '  +0.60%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range('D37').Value = '1.001'
$ws.Range('E37').Value = '  -3.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range('D38').Value = '2.677'
$ws.Range('E38').Value = '  +0.70%  '
$ws.Range('E39').Value = '  -0.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range('D40').Value = '2.729'
$ws.Range('E40').Value = '  +2.11%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '1.142.03'
$ws.Range('E41').Value = '  +7.10%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range('D42').Value = '0.9231'
$ws.Range('E42').Value = '  -2.88%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range('D43').Value = '0.4284'
$ws.Range('E43').Value = '  -0.43%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range('D44').Value = '5.881'
$ws.Range('E44').Value = '  -2.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range('D45').Value = '70.34'
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range('D46').Value = '0.9995'
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range('D47').Value = '102.49'
$ws.Range('E47').Value = '  -0.39%  '
$ws.Range('E48').Value = '  +1.94%  '
$ws.Range('D49').Value = '2.003.95'
$ws.Range('E49').Value = '  +1.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range('D50').Value = '9.210'
$ws.Range('E50').Value = '  +1.05%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range('D51').Value = '6.951'
$ws.Range('E51').Value = '  -1.29%  '
